# Fruta / hortaliza, semanal
# Refresh the weekly "Pera" price records for Terminal Hortofruticola Agro
# Chillan: rows 163-171 roll forward to the next week's figures and the two
# trailing rows that are pushed off the end are appended as new rows 172-173.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 163
$ws.Range("D163").Value = 44610
$ws.Range("K163").Value = "Abate Fettel"
$ws.Range("L163").Value = "Primera"
$ws.Range("M163").Value = 120
$ws.Range("N163").Value = 9000
$ws.Range("O163").Value = 10000
$ws.Range("P163").Value = 9500
$ws.Range("R163").Value = "Región de O'Higgins"
$ws.Range("S163").Value = 594

# Row 164
$ws.Range("D164").Value = 44610
$ws.Range("K164").Value = "Abate Fettel"
$ws.Range("L164").Value = "Segunda"
$ws.Range("M164").Value = 60
$ws.Range("N164").Value = 8000
$ws.Range("O164").Value = 8000
$ws.Range("P164").Value = 8000
$ws.Range("R164").Value = "Región de O'Higgins"
$ws.Range("S164").Value = 500

# Row 165
$ws.Range("D165").Value = 44312
$ws.Range("K165").Value = "Packham's Triumph"
$ws.Range("L165").Value = "Primera"
$ws.Range("M165").Value = 160
$ws.Range("N165").Value = 10000
$ws.Range("O165").Value = 11000
$ws.Range("P165").Value = 10500
$ws.Range("R165").Value = "Provincia de Curicó"
$ws.Range("S165").Value = 656

# Row 166
$ws.Range("D166").Value = 44312
$ws.Range("K166").Value = "Packham's Triumph"
$ws.Range("L166").Value = "Segunda"
$ws.Range("M166").Value = 80
$ws.Range("N166").Value = 9000
$ws.Range("O166").Value = 9000
$ws.Range("P166").Value = 9000
$ws.Range("R166").Value = "Provincia de Curicó"
$ws.Range("S166").Value = 562

# Row 167
$ws.Range("D167").Value = 44399
$ws.Range("K167").Value = "Packham's Triumph"
$ws.Range("L167").Value = "Primera"
$ws.Range("M167").Value = 120
$ws.Range("N167").Value = 9500
$ws.Range("O167").Value = 10000
$ws.Range("P167").Value = 9750
$ws.Range("R167").Value = "Provincia de Curicó"
$ws.Range("S167").Value = 609

# Row 168
$ws.Range("D168").Value = 44399
$ws.Range("K168").Value = "Packham's Triumph"
$ws.Range("L168").Value = "Segunda"
$ws.Range("M168").Value = 120
$ws.Range("N168").Value = 8000
$ws.Range("O168").Value = 8500
$ws.Range("P168").Value = 8250
$ws.Range("R168").Value = "Provincia de Curicó"
$ws.Range("S168").Value = 516

# Row 169
$ws.Range("D169").Value = 44522
$ws.Range("K169").Value = "Packham's Triumph"
$ws.Range("L169").Value = "Primera"
$ws.Range("M169").Value = 120
$ws.Range("N169").Value = 10000
$ws.Range("O169").Value = 11000
$ws.Range("P169").Value = 10500
$ws.Range("R169").Value = "Provincia de Curicó"
$ws.Range("S169").Value = 656

# Row 170
$ws.Range("D170").Value = 44390
$ws.Range("K170").Value = "Packham's Triumph"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 120
$ws.Range("N170").Value = 9000
$ws.Range("O170").Value = 10000
$ws.Range("P170").Value = 9500
$ws.Range("R170").Value = "Provincia de Curicó"
$ws.Range("S170").Value = 594

# Row 171
$ws.Range("D171").Value = 44285
$ws.Range("K171").Value = "Packham's Triumph"
$ws.Range("L171").Value = "Primera"
$ws.Range("M171").Value = 120
$ws.Range("N171").Value = 9000
$ws.Range("O171").Value = 10000
$ws.Range("P171").Value = 9500
$ws.Range("R171").Value = "Provincia de Curicó"
$ws.Range("S171").Value = 594

# Row 172 (new)
$ws.Range("A172").Value = 7
$ws.Range("B172").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C172").Value = "Ñuble"
$ws.Range("D172").Value = 44418
$ws.Range("D172").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E172").Value = 16
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100104
$ws.Range("H172").Value = "Frutos de pepita"
$ws.Range("I172").Value = 100104005
$ws.Range("J172").Value = "Pera"
$ws.Range("K172").Value = "Packham's Triumph"
$ws.Range("L172").Value = "Especial"
$ws.Range("M172").Value = 60
$ws.Range("N172").Value = 9500
$ws.Range("O172").Value = 10000
$ws.Range("P172").Value = 9750
$ws.Range("Q172").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R172").Value = "Provincia de Curicó"
$ws.Range("S172").Value = 609
$ws.Range("T172").Value = 16

# Row 173 (new)
$ws.Range("A173").Value = 7
$ws.Range("B173").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C173").Value = "Ñuble"
$ws.Range("D173").Value = 44418
$ws.Range("D173").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E173").Value = 16
$ws.Range("F173").Value = "Fruta"
$ws.Range("G173").Value = 100104
$ws.Range("H173").Value = "Frutos de pepita"
$ws.Range("I173").Value = 100104005
$ws.Range("J173").Value = "Pera"
$ws.Range("K173").Value = "Packham's Triumph"
$ws.Range("L173").Value = "Primera"
$ws.Range("M173").Value = 60
$ws.Range("N173").Value = 8500
$ws.Range("O173").Value = 9000
$ws.Range("P173").Value = 8750
$ws.Range("Q173").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R173").Value = "Provincia de Curicó"
$ws.Range("S173").Value = 547
$ws.Range("T173").Value = 16
